$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''304.21'
$ws.Range("E2").Value = '''5.67%'

# Row 3
$ws.Range("D3").Value = '''35.12'
$ws.Range("E3").Value = '''13.20%'

# Row 4
$ws.Range("D4").Value = '''5.163'
$ws.Range("E4").Value = '''4.95%'

# Row 5
$ws.Range("D5").Value = '''0.07794'
$ws.Range("E5").Value = '''6.16%'

# Row 6
$ws.Range("D6").Value = '''2.376'
$ws.Range("E6").Value = '''5.47%'

# Row 7
$ws.Range("D7").Value = '''8.057'
$ws.Range("E7").Value = '''4.23%'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''3.971'
$ws.Range("E8").Value = '''6.45%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9241'
$ws.Range("E9").Value = '''1.95%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.09978'
$ws.Range("E10").Value = '''7.84%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1826'
$ws.Range("E11").Value = '''8.48%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.08723'
$ws.Range("E12").Value = '''5.56%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03319'
$ws.Range("E13").Value = '''6.33%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09923'
$ws.Range("E14").Value = '''-0.16%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001483'
$ws.Range("E15").Value = '''-1.00%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005698'
$ws.Range("E16").Value = '''-1.46%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.465'
$ws.Range("E17").Value = '''-1.48%'

# Row 18
$ws.Range("D18").Value = '''2.167'
$ws.Range("E18").Value = '''4.65%'

# Row 19
$ws.Range("D19").Value = '''0.3376'
$ws.Range("E19").Value = '''1.26%'

# Row 20
$ws.Range("D20").Value = '''0.1324'
$ws.Range("E20").Value = '''2.91%'

# Row 21
$ws.Range("D21").Value = '''4.332'
$ws.Range("E21").Value = '''3.15%'

# Row 22
$ws.Range("D22").Value = '''0.2383'
$ws.Range("E22").Value = '''13.36%'

# Row 23
$ws.Range("D23").Value = '''0.04567'
$ws.Range("E23").Value = '''0.63%'

# Row 24
$ws.Range("E24").Value = '''0.76%'

# Row 25
$ws.Range("D25").Value = '''0.004445'
$ws.Range("E25").Value = '''6.98%'

# Row 26
$ws.Range("E26").Value = '''-0.27%'

# Row 27
$ws.Range("D27").Value = '''0.0003697'
$ws.Range("E27").Value = '''8.88%'

# Row 39
$ws.Range("D39").Value = '''0.01783'
$ws.Range("E39").Value = '''13.55%'

# Row 40
$ws.Range("D40").Value = '''0.04799'
$ws.Range("E40").Value = '''7.91%'

# Row 41
$ws.Range("D41").Value = '''0.007756'
$ws.Range("E41").Value = '''4.57%'

# Row 42
$ws.Range("E42").Value = '''6.37%'

# Row 43
$ws.Range("D43").Value = '''0.007160'
$ws.Range("E43").Value = '''-23.04%'

# Row 44
$ws.Range("D44").Value = '''0.002238'
$ws.Range("E44").Value = '''-0.27%'

# Row 45
$ws.Range("D45").Value = '''0.009552'
$ws.Range("E45").Value = '''4.81%'

# Row 46
$ws.Range("E46").Value = '''-0.11%'

# Row 47
$ws.Range("E47").Value = '''-0.11%'

# Row 48
$ws.Range("E48").Value = '''21.85%'

# Row 49
$ws.Range("E49").Value = '''-0.06%'

# Row 50
$ws.Range("E50").Value = '''-0.11%'

# Row 51
$ws.Range("E51").Value = '''-0.11%'
